# Lab 2.xlsx update: add a third data table (rows 61-66) mirroring the
# existing R / V(theo) / P(theo) / V(sim) / P(sim) table, plus a new
# "I_L" column, and move the sheet view/selection down to the new table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 61), columns D:I ----
$ws.Range("D61").Value = "R"
$ws.Range("E61").Value = "V(theo)"
$ws.Range("F61").Value = "P(theo)"
$ws.Range("G61").Value = "V(sim)"
$ws.Range("H61").Value = "P(sim)"
$ws.Range("I61").Value = "I_L"

# Match the formatting of the existing header rows (row 5 / row 29):
# bold+italic font, centered, bottom border, slightly taller row.
$ws.Range("D61:I61").Font.Bold = $true
$ws.Range("D61:I61").Font.Italic = $true
$ws.Range("D61:I61").HorizontalAlignment = -4108
$ws.Range("D61:I61").VerticalAlignment = -4108
$ws.Range("D61:I61").Borders.Item(9).LineStyle = 1
$ws.Range("D61:I61").Borders.Item(9).Weight = -4138
$ws.Rows.Item(61).RowHeight = 15

# ---- Data rows 62-66 ----
$data = @(
    @(2, 5.366, 14.39),
    @(4, 7.589, 14.39),
    @(6, 9.3, 14.39),
    @(8, 10, 13.42),
    @(10, 10, 12)
)

$r = 62
foreach ($row in $data) {
    $ws.Range("D$r").Value = $row[0]
    $ws.Range("E$r").Value = $row[1]
    $ws.Range("F$r").Value = $row[2]
    $ws.Range("G$r").Value = ""
    $ws.Range("H$r").Value = ""
    $ws.Range("I$r").Value = ""
    $r++
}

$ws.Range("D62:I66").HorizontalAlignment = -4108

# ---- Sheet view: scroll to the new table and move the selection ----
$ws.Application.ActiveWindow.ScrollRow = 50
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F67").Select()
